$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.295.59"
$ws.Range("E2").Value = "  +0.81%  "
$ws.Range("D3").Value = "3.487.29"
$ws.Range("E3").Value = "  -0.12%  "
$ws.Range("E4").Value = "  +0.24%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.13"
$ws.Range("E5").Value = "  +0.36%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.28"
$ws.Range("E6").Value = "  +2.12%  "
$ws.Range("D7").Value = "3.487.41"
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  -0.74%  "
$ws.Range("E10").Value = "  -0.38%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.20"
$ws.Range("E11").Value = "  +1.76%  "
$ws.Range("E12").Value = "  -2.59%  "
$ws.Range("D13").Value = "4.079.84"
$ws.Range("E13").Value = "  +0.50%  "
$ws.Range("E14").Value = "  +2.22%  "
$ws.Range("E15").Value = "  +0.67%  "
$ws.Range("D16").Value = "3.485.27"
$ws.Range("E16").Value = "  +0.36%  "
$ws.Range("D17").Value = "64.331.98"
$ws.Range("E17").Value = "  +0.74%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "25.18"
$ws.Range("E18").Value = "  -9.42%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "9.99"
$ws.Range("E19").Value = "  +0.16%  "
$ws.Range("E20").Value = "  +0.70%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.67"
$ws.Range("E21").Value = "  -4.60%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "384.09"
$ws.Range("E22").Value = "  -2.16%  "
$ws.Range("E23").Value = "  -2.03%  "
$ws.Range("D24").Value = "3.626.34"
$ws.Range("E24").Value = "  +0.33%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "74.20"
$ws.Range("E25").Value = "  +2.04%  "
$ws.Range("E26").Value = "  -0.13%  "
$ws.Range("E27").Value = "  -0.78%  "
$ws.Range("E28").Value = "  +1.66%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.54"
$ws.Range("E29").Value = "  -1.12%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  -0.05%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.41"
$ws.Range("E31").Value = "  -0.25%  "
$ws.Range("E32").Value = "  -0.55%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "8.22"
$ws.Range("E33").Value = "  +0.48%  "
$ws.Range("D34").Value = "3.510.32"
$ws.Range("E34").Value = "  +0.91%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.147"
$ws.Range("E36").Value = "  +2.19%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "23.38"
$ws.Range("E37").Value = "  -1.54%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.25"
$ws.Range("E38").Value = "  -1.08%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.84"
$ws.Range("E39").Value = "  -1.38%  "
$ws.Range("E40").Value = "  -2.14%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "161.98"
$ws.Range("E41").Value = "  -4.05%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0776"
$ws.Range("E42").Value = "  -3.77%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.805"
$ws.Range("E43").Value = "  -0.52%  "
$ws.Range("E44").Value = "  +0.22%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "25.44"
$ws.Range("E45").Value = "  -2.35%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "41.73"
$ws.Range("E46").Value = "  -0.10%  "
$ws.Range("E47").Value = "  +0.78%  "
$ws.Range("E48").Value = "  +0.08%  "
$ws.Range("E49").Value = "  +0.42%  "
$ws.Range("D50").Value = "2.463.53"
$ws.Range("E50").Value = "  +1.61%  "
$ws.Range("E51").Value = "  -2.28%  "
